# Auto-generated edit script applying the crypto price/volume update
# described by the commit "Updated cryptos list on Sun Oct  8 13:34:34 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.839.60'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.626.21'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '210.57'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D8').Value = '23.26'
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.857.97'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '1.625.70'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = '0.560'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '65.29'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '27.831.64'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '229.35'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').Value = '0.0₃0721'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = '7.63'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('D21').Value = '0.998'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').Value = '10.10'
$ws.Range('E23').Value = '  -3.37%  '
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('D25').Value = '153.93'
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').Value = '15.50'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('E29').Value = '  +0.44%  '
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = '3.42'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').Value = '3.08'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').Value = '1.394.44'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('E36').Value = '  +10.41%  '
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').Value = '0.0170'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('D46').Value = '1.765.35'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').Value = '87.80'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').Value = '7.60'
$ws.Range('E51').Value = '  +0.50%  '
